$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 224, shifting existing rows down
$ws.Rows(224).Insert()

# Populate the new row 224 with the new weekly data entry
$ws.Range("A224").Value = 10
$ws.Range("B224").Value = "Vega Modelo de Temuco"
$ws.Range("C224").Value = "La Araucanía"
$ws.Range("D224").Value = 44488
$ws.Range("E224").Value = 9
$ws.Range("F224").Value = 100112028
$ws.Range("G224").Value = "Sandia"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Tercera"
$ws.Range("J224").Value = 650
$ws.Range("K224").Value = 900
$ws.Range("L224").Value = 900
$ws.Range("M224").Value = 900
$ws.Range("N224").Value = "$/kilo (volumen en unidades)"
$ws.Range("O224").Value = "Perú"
$ws.Range("P224").Value = 900
$ws.Range("Q224").Value = 1
$ws.Range("R224").Value = "Hortaliza"
